$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing hyperlink in A2 to point to the new URL.
$ws.Range("A2").Value = "https://www.wineenthusiast.com/"
$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://www.wineenthusiast.com/")
$ws.Range("A2").Style = "Hyperlink"

# Fill in the new rows with URLs + hyperlinks (smart-wait-collected links).
$urls = @(
    @{ Cell = "A3"; Url = "https://www.loopnet.com/" },
    @{ Cell = "A4"; Url = "https://www.hannaandersson.com/" },
    @{ Cell = "A5"; Url = "https://www.foco.com/" },
    @{ Cell = "A6"; Url = "https://www.designrush.com/" }
)

foreach ($item in $urls) {
    $cell = $ws.Range($item.Cell)
    $cell.Value = $item.Url
    $ws.Hyperlinks.Add($cell, $item.Url)
    $cell.Style = "Hyperlink"
}

# Match the saved selection state from the source workbook.
$ws.Range("F5").Select() | Out-Null
